# Update "想去人数" (number of people interested) values in both the
# "展览" (Exhibition) and "全部类型" (All types) sheets, which carry the
# same data.
$wb = $excel.ActiveWorkbook

# Map of row -> new F-column value to apply.
$updates = @{
    2  = 8423
    3  = 8015
    11 = 240
    12 = 722
    14 = 2076
    16 = 63
    20 = 48
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
